# LOT2057.xlsx edit: "Docentes responsaveis" names move onto their own rows (13/14),
# stray leftover values in the old A13/A14/B19/C19 cells are dropped, three new Portuguese
# paragraphs are added (Programa resumido / Programa / Bibliografia detail rows), and the
# sheet grows from 23 to 25 rows with the row heights re-tuned to fit the new copy.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop cells that no longer carry data in the new layout ---
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()
$ws.Range("B19").Clear()
$ws.Range("C19").Clear()

# --- Seed formatting for brand-new cells by copying the format from a same-column sibling ---
$ws.Range("B3").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Write the cell values that actually change (unchanged cells are left untouched) ---
$ws.Range("B10").Value = 'Possibilitar aos alunos a execução do projeto proposto e aprovado na disciplina Trabalho de Conclusão de Curso I.'
$ws.Range("C10").Value = 'Possibilitar aos alunos a execução do projeto proposto e aprovado na disciplina Trabalho de Conclusão de Curso I.'
$ws.Range("B13").Value = '1304060 - Maria das Graças de Almeida Felipe'
$ws.Range("C13").Value = '1304060 - Maria das Graças de Almeida Felipe'
$ws.Range("B14").Value = '8853480 - Tatiane da Franca Silva'
$ws.Range("C14").Value = '8853480 - Tatiane da Franca Silva'
$ws.Range("A15").Value = 'Programa resumido:'
$ws.Range("B15").Value = 'Desenvolvimento do trabalho de conclusão de curso, sob orientação de um professor orientador, o qual deve abordar temas relacionados à área de engenharia bioquímica.'
$ws.Range("C15").Value = 'Desenvolvimento do trabalho de conclusão de curso, sob orientação de um professor orientador, o qual deve abordar temas relacionados à área de engenharia bioquímica.'
$ws.Range("A16").Value = 'Short syllabus:'
$ws.Range("B16").Value = 'Development of the course conclusion work under orientation of a leader professor, which must deal with subjects related to Biochemical Engineering.'
$ws.Range("C16").Value = 'Development of the course conclusion work under orientation of a leader professor, which must deal with subjects related to Biochemical Engineering.'
$ws.Range("A17").Value = 'Programa:'
$ws.Range("B17").Value = 'Elaboração de uma monografia de conclusão de curso que apresente: (1) o tema e sua importância, (2) os objetivos, (3) a revisão bibliográfica, (4) a metodologia científica (5) o desenvolvimento do projeto, (6) a análise e a discussão dos resultados, (7) as conclusões e recomendações para trabalhos futuros e (8) referências. O documento deverá atender às normas da ABNT.'
$ws.Range("C17").Value = 'Elaboração de uma monografia de conclusão de curso que apresente: (1) o tema e sua importância, (2) os objetivos, (3) a revisão bibliográfica, (4) a metodologia científica (5) o desenvolvimento do projeto, (6) a análise e a discussão dos resultados, (7) as conclusões e recomendações para trabalhos futuros e (8) referências. O documento deverá atender às normas da ABNT.'
$ws.Range("A18").Value = 'Syllabus:'
$ws.Range("B18").Value = 'Elaboration of a monograph of course conclusion presenting: (1) the subject and its importance, (2) the objectives, (3) the bibliographic revision, (4) the scientific methodology, (5) the development of the project, (6) the analysis and discussion of the results, (7) the conclusion and recommendations for the future works and (8) references. The document must attend to the Brazilian Association of Technical Norms.'
$ws.Range("C18").Value = 'Elaboration of a monograph of course conclusion presenting: (1) the subject and its importance, (2) the objectives, (3) the bibliographic revision, (4) the scientific methodology, (5) the development of the project, (6) the analysis and discussion of the results, (7) the conclusion and recommendations for the future works and (8) references. The document must attend to the Brazilian Association of Technical Norms.'
$ws.Range("A19").Value = 'Avaliação:'
$ws.Range("A20").Value = 'Método:'
$ws.Range("B20").Value = 'Análise da monografia e defesa do trabalho perante uma banca de 3 examinadores, conforme Norma para Trabalho de Conclusão de Curso do curso de Engenharia Bioquímica.'
$ws.Range("C20").Value = 'Análise da monografia e defesa do trabalho perante uma banca de 3 examinadores, conforme Norma para Trabalho de Conclusão de Curso do curso de Engenharia Bioquímica.'
$ws.Range("A21").Value = 'Critério:'
$ws.Range("B21").Value = 'A nota da disciplina será decidida pelos docentes da banca'
$ws.Range("C21").Value = 'A nota da disciplina será decidida pelos docentes da banca'
$ws.Range("A22").Value = 'Norma de recuperação:'
$ws.Range("B22").Value = 'Reapresentação do trabalho modificado para nova avaliação.'
$ws.Range("C22").Value = 'Reapresentação do trabalho modificado para nova avaliação.'
$ws.Range("A23").Value = 'Bibliografia:'
$ws.Range("B23").Value = 'Recomendada pelo Orientador.'
$ws.Range("C23").Value = 'Recomendada pelo Orientador.'
$ws.Range("A24").Value = 'Requisitos:'
$ws.Range("B25").Value = 'LOT2056 -  Trabalho de Conclusão de Curso I  (Requisito)
'
$ws.Range("C25").Value = 'LOT2056 -  Trabalho de Conclusão de Curso I  (Requisito)
'

# --- Row heights: explicit heights where the new layout wants one, AutoFit (-> default) otherwise ---
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 120
$ws.Rows.Item(19).AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 120
$ws.Rows.Item(25).RowHeight = 30
